# fix: handling errors in go-s
#
# Several rows in the "Digital Certificate" column (J) contain a broken
# Google Drive link (https://drive.google.com/file/d//view) where the
# file id never got filled in. For the affected certificate rows, clear
# that bad link out so the cell is blank instead of a dead URL.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose column-J "Digital Certificate" link is broken and needs to
# be cleared.
$rowsToFix = @(6, 8, 15, 17, 22, 33, 35, 40, 48, 55, 63, 65, 68, 69, 75, 82, 83, 88, 93, 94, 100, 106, 107, 108)

foreach ($r in $rowsToFix) {
    $cell = $ws.Cells.Item($r, 10)  # column J = 10
    $cell.Value = ""
}
